$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 540/541 (shifting existing rows 540.. down by 2)
$ws.Rows("540:541").Insert()

# New row 540 data (week of 45131)
$ws.Cells.Item(540, 1).Value = 8
$ws.Cells.Item(540, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(540, 3).Value = "Coquimbo"
$ws.Cells.Item(540, 4).Value = 45131
$ws.Cells.Item(540, 5).Value = 4
$ws.Cells.Item(540, 6).Value = "Fruta"
$ws.Cells.Item(540, 7).Value = 100101
$ws.Cells.Item(540, 8).Value = "Berries"
$ws.Cells.Item(540, 9).Value = 100101007
$ws.Cells.Item(540, 10).Value = "Kiwi"
$ws.Cells.Item(540, 11).Value = "Hayward"
$ws.Cells.Item(540, 12).Value = "Primera"
$ws.Cells.Item(540, 13).Value = 18
$ws.Cells.Item(540, 14).Value = 300000
$ws.Cells.Item(540, 15).Value = 310000
$ws.Cells.Item(540, 16).Value = 305000
$ws.Cells.Item(540, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(540, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(540, 19).Value = 678
$ws.Cells.Item(540, 20).Value = 450

# New row 541 data (week of 45131)
$ws.Cells.Item(541, 1).Value = 8
$ws.Cells.Item(541, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(541, 3).Value = "Coquimbo"
$ws.Cells.Item(541, 4).Value = 45131
$ws.Cells.Item(541, 5).Value = 4
$ws.Cells.Item(541, 6).Value = "Fruta"
$ws.Cells.Item(541, 7).Value = 100101
$ws.Cells.Item(541, 8).Value = "Berries"
$ws.Cells.Item(541, 9).Value = 100101007
$ws.Cells.Item(541, 10).Value = "Kiwi"
$ws.Cells.Item(541, 11).Value = "Hayward"
$ws.Cells.Item(541, 12).Value = "Segunda"
$ws.Cells.Item(541, 13).Value = 10
$ws.Cells.Item(541, 14).Value = 250000
$ws.Cells.Item(541, 15).Value = 260000
$ws.Cells.Item(541, 16).Value = 255000
$ws.Cells.Item(541, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(541, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(541, 19).Value = 567
$ws.Cells.Item(541, 20).Value = 450
